# Regenerate orders with updated distance/size labels.
# Mapping applied throughout all text cells of the sheet:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# (S25 / S20 / NULL are left unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$firstRow = $used.Row
$firstCol = $used.Column

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($firstRow + $r, $firstCol + $c)
        $val = $cell.Value()
        if ($val -is [string]) {
            $newVal = $val -replace 'D64','D69'
            $newVal = $newVal -replace 'D80','D86'
            $newVal = $newVal -replace 'D51','D55'
            $newVal = $newVal -replace 'S30','S31'
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
